$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Model_Home_win (C) and Model_home_win_probability (D) columns
# for rows 2-49 per the updated prediction values.

$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0.4598758518695831
$ws.Cells.Item(3, 4).Value = 0.28061443567276
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0.4257664084434509
$ws.Cells.Item(5, 4).Value = 0.4928802251815796
$ws.Cells.Item(6, 4).Value = 0.2253124266862869
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0.4928478896617889
$ws.Cells.Item(8, 4).Value = 0.504837691783905
$ws.Cells.Item(9, 4).Value = 0.4815551042556763
$ws.Cells.Item(10, 4).Value = 0.2523618638515472
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.4697722196578979
$ws.Cells.Item(12, 4).Value = 0.504837691783905
$ws.Cells.Item(13, 4).Value = 0.3603732287883759
$ws.Cells.Item(14, 4).Value = 0.504837691783905
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 0.504837691783905
$ws.Cells.Item(16, 4).Value = 0.3452500998973846
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 0.504837691783905
$ws.Cells.Item(18, 4).Value = 0.3106700479984283
$ws.Cells.Item(19, 4).Value = 0.2074257433414459
$ws.Cells.Item(20, 4).Value = 0.3105953931808472
$ws.Cells.Item(21, 4).Value = 0.3164584934711456
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0.4747440218925476
$ws.Cells.Item(23, 4).Value = 0.3409002721309662
$ws.Cells.Item(24, 4).Value = 0.504837691783905
$ws.Cells.Item(25, 4).Value = 0.4446874260902405
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 4).Value = 0.504837691783905
$ws.Cells.Item(27, 4).Value = 0.4245058596134186
$ws.Cells.Item(28, 4).Value = 0.1968528181314468
$ws.Cells.Item(29, 4).Value = 0.4138520359992981
$ws.Cells.Item(30, 4).Value = 0.4233212769031525
$ws.Cells.Item(31, 4).Value = 0.4211584627628326
$ws.Cells.Item(32, 4).Value = 0.504837691783905
$ws.Cells.Item(33, 4).Value = 0.2359626442193985
$ws.Cells.Item(34, 4).Value = 0.2342540472745895
$ws.Cells.Item(35, 4).Value = 0.504837691783905
$ws.Cells.Item(36, 4).Value = 0.391911655664444
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 0.494405210018158
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 0.4463329613208771
$ws.Cells.Item(39, 4).Value = 0.392715722322464
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 0.2559272050857544
$ws.Cells.Item(41, 4).Value = 0.504837691783905
$ws.Cells.Item(42, 4).Value = 0.504837691783905
$ws.Cells.Item(43, 4).Value = 0.504837691783905
$ws.Cells.Item(44, 4).Value = 0.504837691783905
$ws.Cells.Item(45, 4).Value = 0.3382506966590881
$ws.Cells.Item(46, 4).Value = 0.504837691783905
$ws.Cells.Item(47, 4).Value = 0.4303871393203735
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 0.4572837054729462
$ws.Cells.Item(49, 4).Value = 0.504837691783905
